$d = $word.ActiveDocument

# --- Fix the capitalisation of the first word ("wejfwkfenwk" -> "Wejfwkfenwk") ---
$d.Content.Find.Execute("wejfwkfenwk", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Wejfwkfenwk", 2)

# --- Append two new paragraphs after the existing text: a blank line, then the sign-off ---
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter()

# Insert the sign-off text plus a throw-away placeholder character. The
# placeholder gives us a safe (non paragraph-boundary) offset to park the
# _GoBack bookmark at before we trim it back off again.
$r = $d.Content
$r.Collapse(0)
$r.InsertAfter("Much love <3X")

# --- Move the _GoBack bookmark from the first paragraph to the end of the
#     newly typed sign-off line (this is where Word leaves it after the user's
#     last real edit). ---
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

$bookmarkPos = $d.Content.End - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

# Remove the placeholder character now that the bookmark is anchored safely.
$placeholder = $d.Range($d.Content.End - 2, $d.Content.End - 1)
$placeholder.Delete()
